$d = $word.ActiveDocument

# Locate the paragraph that currently reads "Si y una polla" and grab a
# Range over its text (excluding the trailing paragraph mark) so that
# InsertXML replaces the run content in place while keeping the
# paragraph's own identity/properties untouched.
$target = $d.Content.Find
$target.Text = "Si y una polla"
$target.Forward = $true
$target.Wrap = 0
$null = $target.Execute()

# After Execute(), Find.Parent is the Range collapsed/extended onto the
# matched text (start/end of the match, paragraph mark excluded) - exactly
# the span whose runs we want to replace.
$r = $target.Parent

$openXmlNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

$pkg = @"
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document $openXmlNs>
<w:body>
<w:p>
<w:r><w:t>-</w:t></w:r>
<w:r><w:t>S</w:t></w:r>
<w:r><w:t>&#237;,</w:t></w:r>
<w:r><w:t xml:space="preserve"> y una polla</w:t></w:r>
<w:r><w:t xml:space="preserve">- dijo Teo y se fue por la </w:t></w:r>
<w:r><w:rPr><w:u w:val="single"/></w:rPr><w:t>puerta</w:t></w:r>
<w:r><w:t>.</w:t></w:r>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
"@

$r.InsertXML($pkg)
